# The edit replaces the tail of the document (the paragraph that begins
# "Isolate is nested within site..." through to the final paragraph
# "Variation due to the different experiments") with a much larger block
# of new paragraphs: the old "Variation due to..." paragraph is pulled up
# right after the (now emptied) "Isolate is nested..." paragraph, a new
# "We describe such situations..." paragraph follows it, and then a long
# run of new paragraphs (an email asking for mixed-model / lmer() advice)
# is appended, ending with the relocated "_GoBack" bookmark and a final
# blank paragraph before the section break.
#
# We rebuild this whole region in one shot via Range.InsertXML so that
# run-level formatting (rPr), <w:proofErr> spell-check markers,
# <w:lastRenderedPageBreak/>, and the bookmark all land exactly where the
# target content expects them.

$d = $word.ActiveDocument

# Find the paragraph that currently starts the region to be rebuilt.
$startText = "Isolate is nested within site and sampling date is nested within isolate"
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "$startText*") {
        $anchor = $d.Paragraphs.Item($i)
        break
    }
}
if ($anchor -eq $null) {
    throw "Could not locate the 'Isolate is nested...' paragraph"
}

# The region runs from the start of that paragraph through the end of the
# document body's story (the last paragraph, "Variation due to the
# different experiments", right before the final section break).
$region = $d.Range($anchor.Range.Start, $d.Content.End)

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$p01 = '<w:p ' + $wns + '><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p>'

$p02 = '<w:p ' + $wns + '><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Variation due to the different experiments</w:t></w:r></w:p>'

$p03 = '<w:p ' + $wns + '><w:r><w:t>We describe such situations as having partially crossed grouping factors for the random effects.</w:t></w:r></w:p>'

$p04 = '<w:p ' + $wns + '><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p>'

$p05 = '<w:p ' + $wns + '><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>' +
       '<w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Completely crossed each isolate in each vegetation type for each sampling date</w:t></w:r></w:p>'

$p06 = '<w:p ' + $wns + '>' +
       '<w:r><w:t xml:space="preserve">At this point we will fit models that have random effects for </w:t></w:r>' +
       '<w:r><w:t>isolate</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">, </w:t></w:r>' +
       '<w:r><w:t>incubation time</w:t></w:r>' +
       '<w:r><w:t xml:space="preserve">, and </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>experiemnt</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> (or the </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>dept:service</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> combination) to these data. In the next </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>ch</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '</w:p>'

$p07 = '<w:p ' + $wns + '><w:r><w:rPr><w:rFonts w:ascii="Merriweather" w:hAnsi="Merriweather"/><w:color w:val="5C5C5C"/>' +
       '<w:sz w:val="21"/><w:szCs w:val="21"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr>' +
       '<w:t>explicitly in variation among and by groups. This is where a mixed-effect modeling framework is useful</w:t></w:r></w:p>'

$p08 = '<w:p ' + $wns + '/>'

$p09 = '<w:p ' + $wns + '><w:r><w:t>standing of how to explain the relationships among the fixed and random effects in terms of the levels of the hierarchy.</w:t></w:r></w:p>'

$p10 = '<w:p ' + $wns + '/>'

$p11 = '<w:p ' + $wns + '><w:r><w:t>Hello,</w:t></w:r></w:p>'

$p12 = '<w:p ' + $wns + '/>'

$p13 = '<w:p ' + $wns + '>' +
       '<w:r><w:t xml:space="preserve">I have a dataset which combines 3 experiments measuring the decomposition of fungi. ' +
       'The experiments were conducted at three different sites, differing their dominant vegetation ' +
       '(i.e., prairie, oak savanna and forest). The species of fungi decomposed and the times the fungi ' +
       'were decomposed differ among the 3 experiments/sites. There are two species of fungi which were ' +
       'included at all the sites, but not all combinations of factors occur (thus my design is partially </w:t></w:r>' +
       '<w:r><w:lastRenderedPageBreak/><w:t>crossed). I have decided to use a mixed linear effects model to deal with any variation caused by differences among experiments-</w:t></w:r>' +
       '<w:bookmarkStart w:id="2" w:name="_GoBack"/>' +
       '<w:r><w:t xml:space="preserve">treating the partially crossed grouping factors as a random effect. </w:t></w:r>' +
       '<w:bookmarkEnd w:id="2"/>' +
       '</w:p>'

$p14 = '<w:p ' + $wns + '/>'

$p15 = '<w:p ' + $wns + '>' +
       '<w:r><w:t xml:space="preserve">I perform my statistical analyses in R and I am looking for someone to help me confirm that I am using the correct </w:t></w:r>' +
       '<w:proofErr w:type="spellStart"/><w:r><w:t>lmer</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
       '<w:r><w:t xml:space="preserve"> () model syntax to specify relationships among fixed and random effects. </w:t></w:r>' +
       '</w:p>'

$p16 = '<w:p ' + $wns + '><w:pPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p>'

$xml = $p01 + $p02 + $p03 + $p04 + $p05 + $p06 + $p07 + $p08 + $p09 + $p10 +
       $p11 + $p12 + $p13 + $p14 + $p15 + $p16

$region.InsertXML($xml) | Out-Null
